# Auto-generated: apply scheduled-runner market price updates across sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 2006.8462
$ws.Range("I2").Value = 1841.4286
$ws.Range("J2").Value = 2199.8333
$ws.Range("K2").Value = 1841.4286
$ws.Range("L2").Value = 2199.8333
$ws.Range("M2").Value = -1728.4286
$ws.Range("N2").Value = -2425.8333

$ws.Range("H11").Value = 978.5
$ws.Range("I11").Value = 978.5
$ws.Range("K11").Value = 978.5
$ws.Range("M11").Value = -838.5

$ws.Range("H17").Value = 1815
$ws.Range("I17").Value = 1471.75
$ws.Range("K17").Value = 4415.25
$ws.Range("M17").Value = -4247.25

$ws.Range("H18").Value = 907.8
$ws.Range("I18").Value = 907.8
$ws.Range("K18").Value = 907.8
$ws.Range("M18").Value = -623.8

$ws.Range("H39").Value = 620
$ws.Range("I39").Value = 517.5
$ws.Range("J39").Value = 825
$ws.Range("K39").Value = 1552.5
$ws.Range("L39").Value = 2475
$ws.Range("M39").Value = -1256.5
$ws.Range("N39").Value = -3067

$ws.Range("H74").Value = 4901.8
$ws.Range("I74").Value = 4901.8
$ws.Range("K74").Value = 4901.8
$ws.Range("M74").Value = -3965.8

$ws.Range("H77").Value = 4901.8
$ws.Range("I77").Value = 4901.8
$ws.Range("K77").Value = 24509
$ws.Range("M77").Value = -19829

$ws.Range("H92").Value = 1389.5938
$ws.Range("I92").Value = 543.04
$ws.Range("J92").Value = 4413
$ws.Range("K92").Value = 543.04
$ws.Range("L92").Value = 4413
$ws.Range("M92").Value = 704.96
$ws.Range("N92").Value = -6909

$ws.Range("H98").Value = 1729.6923
$ws.Range("I98").Value = 624.125
$ws.Range("K98").Value = 624.125
$ws.Range("M98").Value = 873.875

$ws.Range("H112").Value = 4527.591
$ws.Range("J112").Value = 4797.2197
$ws.Range("L112").Value = 14391.6591
$ws.Range("N112").Value = -16607.6591

$ws.Range("H122").Value = 1729.6923
$ws.Range("I122").Value = 624.125
$ws.Range("K122").Value = 1872.375
$ws.Range("M122").Value = 577.625

$ws.Range("H132").Value = 2816928.2
$ws.Range("I132").Value = 2931787.2
$ws.Range("K132").Value = 8795361.600000001
$ws.Range("M132").Value = -8792831.600000001

$ws.Range("H137").Value = 6672.3193
$ws.Range("I137").Value = 9657.710999999999
$ws.Range("J137").Value = 3335.7058
$ws.Range("K137").Value = 28973.133
$ws.Range("L137").Value = 10007.1174
$ws.Range("M137").Value = -26423.133
$ws.Range("N137").Value = -15107.1174

$ws.Range("H141").Value = 1519.9375
$ws.Range("I141").Value = 1519.9375
$ws.Range("K141").Value = 4559.8125
$ws.Range("M141").Value = 620.1875

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("N3").ClearContents()

$ws.Range("H61").Value = 3267.9033
$ws.Range("I61").Value = 2257.1304
$ws.Range("J61").Value = 6173.875
$ws.Range("K61").Value = 2257.1304
$ws.Range("L61").Value = 6173.875
$ws.Range("M61").Value = -2045.1304
$ws.Range("N61").Value = -6597.875

$ws.Range("H74").Value = 242514.38
$ws.Range("I74").Value = 251125.61
$ws.Range("J74").Value = 1400
$ws.Range("K74").Value = 251125.61
$ws.Range("L74").Value = 1400
$ws.Range("M74").Value = -250251.61
$ws.Range("N74").Value = -3148

$ws.Range("H77").Value = 242514.38
$ws.Range("I77").Value = 251125.61
$ws.Range("J77").Value = 1400
$ws.Range("K77").Value = 1255628.05
$ws.Range("L77").Value = 7000
$ws.Range("M77").Value = -1251260.05
$ws.Range("N77").Value = -15736

$ws.Range("H132").Value = 1304.7567
$ws.Range("I132").Value = 1104.5
$ws.Range("K132").Value = 3313.5
$ws.Range("M132").Value = -783.5

$ws.Range("H136").Value = 3267.9033
$ws.Range("I136").Value = 2257.1304
$ws.Range("J136").Value = 6173.875
$ws.Range("K136").Value = 6771.3912
$ws.Range("L136").Value = 18521.625
$ws.Range("M136").Value = -4221.3912
$ws.Range("N136").Value = -23621.625

$ws.Range("H139").Value = 78369.75
$ws.Range("J139").Value = 78369.75
$ws.Range("L139").Value = 78369.75
$ws.Range("N139").Value = -88649.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1073.4857
$ws.Range("I94").Value = 1157.44
$ws.Range("K94").Value = 1157.44
$ws.Range("M94").Value = -706.4400000000001

$ws.Range("H134").Value = 3005.1353
$ws.Range("I134").Value = 2876.2058
$ws.Range("K134").Value = 8628.617400000001
$ws.Range("M134").Value = -6093.617400000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2225135.8
$ws.Range("I31").Value = 3032899.5
$ws.Range("J31").Value = 3785.5
$ws.Range("K31").Value = 3032899.5
$ws.Range("L31").Value = 3785.5
$ws.Range("M31").Value = -3032604.5
$ws.Range("N31").Value = -4375.5

$ws.Range("H34").Value = 2225135.8
$ws.Range("I34").Value = 3032899.5
$ws.Range("J34").Value = 3785.5
$ws.Range("K34").Value = 3032899.5
$ws.Range("L34").Value = 3785.5
$ws.Range("M34").Value = -3032697.5
$ws.Range("N34").Value = -4189.5

$ws.Range("H58").Value = 1904.65
$ws.Range("I58").Value = 1510.9706
$ws.Range("J58").Value = 4135.5
$ws.Range("K58").Value = 1510.9706
$ws.Range("L58").Value = 4135.5
$ws.Range("M58").Value = -1307.9706
$ws.Range("N58").Value = -4541.5

$ws.Range("H136").Value = 1904.65
$ws.Range("I136").Value = 1510.9706
$ws.Range("J136").Value = 4135.5
$ws.Range("K136").Value = 4532.9118
$ws.Range("L136").Value = 12406.5
$ws.Range("M136").Value = -1982.9118
$ws.Range("N136").Value = -17506.5

$ws.Range("H139").Value = 51830.75
$ws.Range("J139").Value = 51830.75
$ws.Range("L139").Value = 51830.75
$ws.Range("N139").Value = -62110.75

$ws.Range("H140").Value = 119888
$ws.Range("J140").Value = 119888
$ws.Range("L140").Value = 119888
$ws.Range("N140").Value = -130248

$ws.Range("H141").Value = 249999.28
$ws.Range("J141").Value = 249999.28
$ws.Range("L141").Value = 249999.28
$ws.Range("N141").Value = -260359.28

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H103").Value = 446.33334
$ws.Range("I103").Value = 446.33334
$ws.Range("K103").Value = 1339.00002
$ws.Range("M103").Value = -460.0000199999999

$ws.Range("H113").Value = 2388.2273
$ws.Range("J113").Value = 2021.4667
$ws.Range("L113").Value = 6064.4001
$ws.Range("N113").Value = -10404.4001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 5561.125
$ws.Range("J122").Value = 6000
$ws.Range("L122").Value = 18000
$ws.Range("N122").Value = -22900

$ws.Range("H132").Value = 1538.0476
$ws.Range("I132").Value = 1334.8
$ws.Range("K132").Value = 4004.4
$ws.Range("M132").Value = -1474.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H13").Value = 15500
$ws.Range("I13").Value = 26500
$ws.Range("K13").Value = 26500
$ws.Range("M13").Value = -26360

$ws.Range("H17").Value = 20000
$ws.Range("I17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("M17").ClearContents()

$ws.Range("H22").Value = 1046.1072
$ws.Range("J22").Value = 1411.1538
$ws.Range("L22").Value = 1411.1538
$ws.Range("N22").Value = -2001.1538

$ws.Range("H27").Value = 1046.1072
$ws.Range("J27").Value = 1411.1538
$ws.Range("L27").Value = 1411.1538
$ws.Range("N27").Value = -1625.1538

$ws.Range("H103").Value = 31750
$ws.Range("J103").Value = 31750
$ws.Range("L103").Value = 31750
$ws.Range("N103").Value = -34094

$ws.Range("H132").Value = 3722.7778
$ws.Range("I132").Value = 3722.7778
$ws.Range("K132").Value = 11168.3334
$ws.Range("M132").Value = -8638.3334

$ws.Range("H136").Value = 2596.04
$ws.Range("I136").Value = 2404.5789
$ws.Range("J136").Value = 3202.3333
$ws.Range("K136").Value = 7213.736699999999
$ws.Range("L136").Value = 9606.999899999999
$ws.Range("M136").Value = -4663.736699999999
$ws.Range("N136").Value = -14706.9999

$ws.Range("H137").Value = 0
$ws.Range("I137").Value = 0
$ws.Range("K137").Value = 0
$ws.Range("M137").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 68219.14999999999
$ws.Range("J122").Value = 4517.5
$ws.Range("L122").Value = 13552.5
$ws.Range("N122").Value = -18452.5

$ws.Range("H131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").ClearContents()

$ws.Range("H132").Value = 27471.834
$ws.Range("I132").Value = 28281.207
$ws.Range("K132").Value = 84843.621
$ws.Range("M132").Value = -82313.621

$ws.Range("H136").Value = 21438.107
$ws.Range("I136").Value = 22023.361
$ws.Range("J136").Value = 369
$ws.Range("K136").Value = 66070.083
$ws.Range("L136").Value = 1107
$ws.Range("M136").Value = -63520.083
$ws.Range("N136").Value = -6207

Write-Output "Applied scheduled profit updates"